# Update the handback status report timestamps to reflect the latest
# generation run ("Generate Report for Handback").

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 988f0d34 row (row 3)
$overview.Range("G3").Value = "2017-02-09 08:08:04"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for row 3
$zhcn.Range("H3").Value = "2017-02-09 08:07:46"
$zhcn.Range("L3").Value = "2017-02-09 08:08:43"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for row 3
$dede.Range("H3").Value = "2017-02-09 08:08:04"
$dede.Range("L3").Value = "2017-02-09 08:09:08"
